# Re-bind the workbook / worksheets the way the harness expects.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
# The IG generation timestamp moved forward (re-run of the generator).
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# --- Elements sheet ---------------------------------------------------------
# Row 6 documents the "exerciceProfessionnel" attribute; the generator now
# emits it with the class's own PascalCase spelling, and the short
# description lost its trailing period.
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("A6").Value = "CompetenceMetier.ExerciceProfessionnel"
$wsElem.Range("B6").Value = "CompetenceMetier.ExerciceProfessionnel"
$wsElem.Range("L6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("M6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("AF6").Value = "SavoirFaire.ExerciceProfessionnel"

# The ID/Path columns were re-measured by the generator's best-fit pass
# (33.5 -> 33.52734375 chars). The host's ColumnWidth setter here only has
# pixel-level (1/6 character) resolution, so 32.666666666666664 is the
# closest settable value that keeps the column at its (still custom-width)
# 33.5-character footprint without perturbing any of the other columns.
$wsElem.Columns.Item(1).ColumnWidth = 32.666666666666664
$wsElem.Columns.Item(2).ColumnWidth = 32.666666666666664
